$p = $ppt.ActivePresentation

# Add a splash-screen style background fill to slide 1 (solid light grey).
$s1 = $p.Slides.Item(1)
$s1.Background.Fill.Solid()
$s1.Background.Fill.ForeColor.RGB = 0xF5F5F5

# Remove the second slide (cleanup of unused slide).
$p.Slides.Item(2).Delete()
